$wb = $excel.ActiveWorkbook

# --- Sheet "Status" (Worksheets(1)): add two new rows to the Table3 table ---
$wsStatus = $wb.Worksheets.Item(1)
$loStatus = $wsStatus.ListObjects.Item(1)

# Expand the table so it covers the two new data rows (B4:G11 -> B4:G13)
$loStatus.Resize($wsStatus.Range("B4:G13"))

# Row 12: id=8, code=20, status_type="ROOM TRANSACTION", name="Check in", seq_num=1
$wsStatus.Range("B12").Value = 8
$wsStatus.Range("C12").Value = 20
$wsStatus.Range("D12").Value = "ROOM TRANSACTION"
$wsStatus.Range("E12").Value = "Check in"
$wsStatus.Range("F12").Value = 1
$wsStatus.Range("F12").HorizontalAlignment = -4131

# Row 13: id=9, code=21, status_type="ROOM TRANSACTION", name="Check out", seq_num=2
$wsStatus.Range("B13").Value = 9
$wsStatus.Range("C13").Value = 21
$wsStatus.Range("D13").Value = "ROOM TRANSACTION"
$wsStatus.Range("E13").Value = "Check out"
$wsStatus.Range("F13").Value = 2
$wsStatus.Range("F13").HorizontalAlignment = -4131

# --- Sheet "Cfg Utilities" (Worksheets(4)): update selection / scroll position only ---
$wsCfg = $wb.Worksheets.Item(4)
$wsCfg.Range("F13").Select()

# --- Restore "Status" as the active sheet / update its selection ---
$wsStatus.Select()
$wsStatus.Range("E14").Select()
